$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New voltage-magnitude results for the 380 kV case (columns B:F, rows 2-25)
$dataBF = New-Object 'object[,]' 24,5
$dataBF[0,0] = 1.02
$dataBF[0,1] = 1.038260940891049
$dataBF[0,2] = 1.040913239600666
$dataBF[0,3] = 1.052049505942946
$dataBF[0,4] = 1.059980282238356
$dataBF[1,0] = 1.02
$dataBF[1,1] = 1.039100721899552
$dataBF[1,2] = 1.041534614806333
$dataBF[1,3] = 1.052881849971335
$dataBF[1,4] = 1.060899726521548
$dataBF[2,0] = 1.02
$dataBF[2,1] = 1.039644658303895
$dataBF[2,2] = 1.041937061483112
$dataBF[2,3] = 1.053421342129909
$dataBF[2,4] = 1.061495719936049
$dataBF[3,0] = 1.02
$dataBF[3,1] = 1.039873457380676
$dataBF[3,2] = 1.042106338271637
$dataBF[3,3] = 1.053648360717133
$dataBF[3,4] = 1.061746525505026
$dataBF[4,0] = 1.02
$dataBF[4,1] = 1.039911881226495
$dataBF[4,2] = 1.042134765713445
$dataBF[4,3] = 1.053686490757196
$dataBF[4,4] = 1.061788651468125
$dataBF[5,0] = 1.02
$dataBF[5,1] = 1.039647715026478
$dataBF[5,2] = 1.041939323023344
$dataBF[5,3] = 1.053424374716517
$dataBF[5,4] = 1.061499070232494
$dataBF[6,0] = 1.02
$dataBF[6,1] = 1.03854463568322
$dataBF[6,2] = 1.04112315754173
$dataBF[6,3] = 1.052330611047167
$dataBF[6,4] = 1.060290793946174
$dataBF[7,0] = 1.02
$dataBF[7,1] = 1.03660508717151
$dataBF[7,2] = 1.039687917806609
$dataBF[7,3] = 1.050410306031985
$dataBF[7,4] = 1.058169788998531
$dataBF[8,0] = 1.02
$dataBF[8,1] = 1.035314981735325
$dataBF[8,2] = 1.03873317220993
$dataBF[8,3] = 1.049134944517566
$dataBF[8,4] = 1.056761359855687
$dataBF[9,0] = 1.02
$dataBF[9,1] = 1.034757065086245
$dataBF[9,2] = 1.038320271215164
$dataBF[9,3] = 1.048583869273349
$dataBF[9,4] = 1.056152838996284
$dataBF[10,0] = 1.02
$dataBF[10,1] = 1.034549937814919
$dataBF[10,2] = 1.038166979727213
$dataBF[10,3] = 1.048379352144919
$dataBF[10,4] = 1.055927010114089
$dataBF[11,0] = 1.02
$dataBF[11,1] = 1.034594362414561
$dataBF[11,2] = 1.038199857700351
$dataBF[11,3] = 1.048423213736315
$dataBF[11,4] = 1.055975441967154
$dataBF[12,0] = 1.02
$dataBF[12,1] = 1.03473994167324
$dataBF[12,2] = 1.038307598485408
$dataBF[12,3] = 1.048566960207753
$dataBF[12,4] = 1.056134167761738
$dataBF[13,0] = 1.02
$dataBF[13,1] = 1.03482965219098
$dataBF[13,2] = 1.038373991568184
$dataBF[13,3] = 1.048655550627584
$dataBF[13,4] = 1.056231990897154
$dataBF[14,0] = 1.02
$dataBF[14,1] = 1.03535202386186
$dataBF[14,2] = 1.038760585980963
$dataBF[14,3] = 1.04917154225671
$dataBF[14,4] = 1.05680177373025
$dataBF[15,0] = 1.02
$dataBF[15,1] = 1.035679884302572
$dataBF[15,2] = 1.039003224211818
$dataBF[15,3] = 1.049495523247903
$dataBF[15,4] = 1.057159543008046
$dataBF[16,0] = 1.02
$dataBF[16,1] = 1.035871187914039
$dataBF[16,2] = 1.039144800096094
$dataBF[16,3] = 1.049684608137447
$dataBF[16,4] = 1.057368352775961
$dataBF[17,0] = 1.02
$dataBF[17,1] = 1.035936429017409
$dataBF[17,2] = 1.039193082092573
$dataBF[17,3] = 1.049749100185575
$dataBF[17,4] = 1.057439573336726
$dataBF[18,0] = 1.02
$dataBF[18,1] = 1.035644700918715
$dataBF[18,2] = 1.038977186334205
$dataBF[18,3] = 1.049460751529454
$dataBF[18,4] = 1.05712114439423
$dataBF[19,0] = 1.02
$dataBF[19,1] = 1.034697069221653
$dataBF[19,2] = 1.038275869333905
$dataBF[19,3] = 1.048524625561215
$dataBF[19,4] = 1.05608742135698
$dataBF[20,0] = 1.02
$dataBF[20,1] = 1.034101880090533
$dataBF[20,2] = 1.037835377500933
$dataBF[20,3] = 1.047937069530026
$dataBF[20,4] = 1.055438653230021
$dataBF[21,0] = 1.02
$dataBF[21,1] = 1.034417341404887
$dataBF[21,2] = 1.038068846961072
$dataBF[21,3] = 1.048248446490224
$dataBF[21,4] = 1.055782465551
$dataBF[22,0] = 1.02
$dataBF[22,1] = 1.035660598564987
$dataBF[22,2] = 1.038988951578783
$dataBF[22,3] = 1.049476463025645
$dataBF[22,4] = 1.057138494676587
$dataBF[23,0] = 1.02
$dataBF[23,1] = 1.037105997725217
$dataBF[23,2] = 1.040058602129682
$dataBF[23,3] = 1.050905905138649
$dataBF[23,4] = 1.058717144937203
$ws.Range("B2:F25").Value = $dataBF

# New voltage-magnitude results for the 380 kV case (columns I:N, rows 2-25)
$dataIN = New-Object 'object[,]' 24,6
$dataIN[0,0] = 1.042238535863114
$dataIN[0,1] = 1.043359340641833
$dataIN[0,2] = 1.043694241123195
$dataIN[0,3] = 1.054799314032762
$dataIN[0,4] = 1.062708308513658
$dataIN[0,5] = 1.018380496326309
$dataIN[1,0] = 1.042461467614356
$dataIN[1,1] = 1.043844391716764
$dataIN[1,2] = 1.044126634127934
$dataIN[1,3] = 1.055444375993535
$dataIN[1,4] = 1.063441822878299
$dataIN[1,5] = 1.01854258267677
$dataIN[2,0] = 1.04260466508376
$dataIN[2,1] = 1.044158125593912
$dataIN[2,2] = 1.044406104102086
$dataIN[2,3] = 1.05586203394583
$dataIN[2,4] = 1.063916873560452
$dataIN[2,5] = 1.018647386375969
$dataIN[3,0] = 1.042664612139034
$dataIN[3,1] = 1.044289987831182
$dataIN[3,2] = 1.044523516167676
$dataIN[3,3] = 1.056037678086328
$dataIN[3,4] = 1.064116683341717
$dataIN[3,5] = 1.018691426978384
$dataIN[4,0] = 1.042674662645743
$dataIN[4,1] = 1.044312126200593
$dataIN[4,2] = 1.044543225614252
$dataIN[4,3] = 1.056067173030136
$dataIN[4,4] = 1.064150238012546
$dataIN[4,5] = 1.018698820469908
$dataIN[5,0] = 1.042605467094948
$dataIN[5,1] = 1.044159887668766
$dataIN[5,2] = 1.044407673272316
$dataIN[5,3] = 1.055864380674583
$dataIN[5,4] = 1.063919543043541
$dataIN[5,5] = 1.018647974923359
$dataIN[6,0] = 1.042314094528267
$dataIN[6,1] = 1.043523291430563
$dataIN[6,2] = 1.043840435341261
$dataIN[6,3] = 1.055017261096431
$dataIN[6,4] = 1.062956115674516
$dataIN[6,5] = 1.018435289849744
$dataIN[7,0] = 1.041792618115486
$dataIN[7,1] = 1.042400609872633
$dataIN[7,2] = 1.04283851186971
$dataIN[7,3] = 1.053526581164024
$dataIN[7,4] = 1.061261702576988
$dataIN[7,5] = 1.018059941088907
$dataIN[8,0] = 1.041439606514577
$dataIN[8,1] = 1.041651608445245
$dataIN[8,2] = 1.042169033820652
$dataIN[8,3] = 1.052534260685406
$dataIN[8,4] = 1.060134384402748
$dataIN[8,5] = 1.017809351000965
$dataIN[9,0] = 1.04128548665342
$dataIN[9,1] = 1.041327167496552
$dataIN[9,2] = 1.04187879439821
$dataIN[9,3] = 1.052104940758537
$dataIN[9,4] = 1.059646805953446
$dataIN[9,5] = 1.017700763378555
$dataIN[10,0] = 1.041228050594378
$dataIN[10,1] = 1.041206639069088
$dataIN[10,2] = 1.041770934987393
$dataIN[10,3] = 1.051945527963293
$dataIN[10,4] = 1.059465782966216
$dataIN[10,5] = 1.017660417449597
$dataIN[11,0] = 1.041240379368283
$dataIN[11,1] = 1.041232493558874
$dataIN[11,2] = 1.041794073511875
$dataIN[11,3] = 1.051979720003381
$dataIN[11,4] = 1.059504609123753
$dataIN[11,5] = 1.017669072308538
$dataIN[12,0] = 1.041280742823716
$dataIN[12,1] = 1.041317204911261
$dataIN[12,2] = 1.041869879750845
$dataIN[12,3] = 1.052091762501684
$dataIN[12,4] = 1.059631840786534
$dataIN[12,5] = 1.017697428606634
$dataIN[13,0] = 1.04130558704817
$dataIN[13,1] = 1.041369396204957
$dataIN[13,2] = 1.041916579678506
$dataIN[13,3] = 1.052160803012269
$dataIN[13,4] = 1.05971024377067
$dataIN[13,5] = 1.017714898326952
$dataIN[14,0] = 1.041449808385267
$dataIN[14,1] = 1.041673138128309
$dataIN[14,2] = 1.042188288759641
$dataIN[14,3] = 1.05256276097252
$dataIN[14,4] = 1.060166755263211
$dataIN[14,5] = 1.017816555945947
$dataIN[15,0] = 1.041539936869736
$dataIN[15,1] = 1.041863636544205
$dataIN[15,2] = 1.042358631537707
$dataIN[15,3] = 1.05281499635634
$dataIN[15,4] = 1.060453263291352
$dataIN[15,5] = 1.017880301851663
$dataIN[16,0] = 1.041592385329406
$dataIN[16,1] = 1.041974739545801
$dataIN[16,2] = 1.042457955652496
$dataIN[16,3] = 1.052962155759965
$dataIN[16,4] = 1.060620432287751
$dataIN[16,5] = 1.017917475967377
$dataIN[17,0] = 1.04161024818944
$dataIN[17,1] = 1.042012620811664
$dataIN[17,2] = 1.042491816815536
$dataIN[17,3] = 1.053012339167551
$dataIN[17,4] = 1.060677441658311
$dataIN[17,5] = 1.017930150045107
$dataIN[18,0] = 1.041530279545868
$dataIN[18,1] = 1.041843199048214
$dataIN[18,2] = 1.042340358872453
$dataIN[18,3] = 1.052787930279742
$dataIN[18,4] = 1.060422518124574
$dataIN[18,5] = 1.017873463322134
$dataIN[19,0] = 1.041268862005657
$dataIN[19,1] = 1.041292259978546
$dataIN[19,2] = 1.041847558116988
$dataIN[19,3] = 1.052058767228174
$dataIN[19,4] = 1.059594371864161
$dataIN[19,5] = 1.017689078701725
$dataIN[20,0] = 1.041103404832829
$dataIN[20,1] = 1.040945767179893
$dataIN[20,2] = 1.041537416954025
$dataIN[20,3] = 1.051600636694429
$dataIN[20,4] = 1.059074178071653
$dataIN[20,5] = 1.017573081579074
$dataIN[21,0] = 1.041191220220296
$dataIN[21,1] = 1.04112945823245
$dataIN[21,2] = 1.041701856487351
$dataIN[21,3] = 1.051843469238034
$dataIN[21,4] = 1.059349895236533
$dataIN[21,5] = 1.017634580071782
$dataIN[22,0] = 1.041534643651751
$dataIN[22,1] = 1.041852433908606
$dataIN[22,2] = 1.042348615608369
$dataIN[22,3] = 1.052800160167444
$dataIN[22,4] = 1.060436410376257
$dataIN[22,5] = 1.017876553383397
$dataIN[23,0] = 1.041928380456134
$dataIN[23,1] = 1.042690951509183
$dataIN[23,2] = 1.043097807707815
$dataIN[23,3] = 1.053911705526243
$dataIN[23,4] = 1.061699352039373
$dataIN[23,5] = 1.018157042668257
$ws.Range("I2:N25").Value = $dataIN
